$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 8; $row++) {
    $cell = $ws.Range("C$row")
    if ($cell.Value2 -eq 46075) {
        $cell.Value2 = 46076
    }
}
